# Insert 3 new rows above row 1144 (new week of price data), shifting the
# rest of the table down, then populate the 3 new rows with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 1144 (existing row 1144 and below shift down to 1147+)
$ws.Rows.Item(1144).Resize(3).Insert()

# New row 1144: Conconina(o)
$ws.Cells.Item(1144, 1).Value = 10
$ws.Cells.Item(1144, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1144, 3).Value = "La Araucanía"
$ws.Cells.Item(1144, 4).Value = 44753
$ws.Cells.Item(1144, 5).Value = 9
$ws.Cells.Item(1144, 6).Value = 100112033
$ws.Cells.Item(1144, 7).Value = "Lechuga"
$ws.Cells.Item(1144, 8).Value = "Conconina(o)"
$ws.Cells.Item(1144, 9).Value = "Primera"
$ws.Cells.Item(1144, 10).Value = 400
$ws.Cells.Item(1144, 11).Value = 7000
$ws.Cells.Item(1144, 12).Value = 7000
$ws.Cells.Item(1144, 13).Value = 7000
$ws.Cells.Item(1144, 14).Value = "$/caja 10 unidades"
$ws.Cells.Item(1144, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1144, 16).Value = 700
$ws.Cells.Item(1144, 17).Value = 10
$ws.Cells.Item(1144, 18).Value = "Hortaliza"

# New row 1145: Escarola
$ws.Cells.Item(1145, 1).Value = 10
$ws.Cells.Item(1145, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1145, 3).Value = "La Araucanía"
$ws.Cells.Item(1145, 4).Value = 44753
$ws.Cells.Item(1145, 5).Value = 9
$ws.Cells.Item(1145, 6).Value = 100112033
$ws.Cells.Item(1145, 7).Value = "Lechuga"
$ws.Cells.Item(1145, 8).Value = "Escarola"
$ws.Cells.Item(1145, 9).Value = "Primera"
$ws.Cells.Item(1145, 10).Value = 900
$ws.Cells.Item(1145, 11).Value = 12000
$ws.Cells.Item(1145, 12).Value = 12000
$ws.Cells.Item(1145, 13).Value = 12000
$ws.Cells.Item(1145, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(1145, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(1145, 16).Value = 800
$ws.Cells.Item(1145, 17).Value = 15
$ws.Cells.Item(1145, 18).Value = "Hortaliza"

# New row 1146: Marina
$ws.Cells.Item(1146, 1).Value = 10
$ws.Cells.Item(1146, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1146, 3).Value = "La Araucanía"
$ws.Cells.Item(1146, 4).Value = 44753
$ws.Cells.Item(1146, 5).Value = 9
$ws.Cells.Item(1146, 6).Value = 100112033
$ws.Cells.Item(1146, 7).Value = "Lechuga"
$ws.Cells.Item(1146, 8).Value = "Marina"
$ws.Cells.Item(1146, 9).Value = "Primera"
$ws.Cells.Item(1146, 10).Value = 500
$ws.Cells.Item(1146, 11).Value = 8000
$ws.Cells.Item(1146, 12).Value = 8000
$ws.Cells.Item(1146, 13).Value = 8000
$ws.Cells.Item(1146, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(1146, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1146, 16).Value = 533
$ws.Cells.Item(1146, 17).Value = 15
$ws.Cells.Item(1146, 18).Value = "Hortaliza"
